# This script re-orders the data that lives in rows 16, 18, 19, 20 and 21
# of the "Artfynd" sheet. Row 17 is left untouched. The underlying records
# (identified by the tuple of columns A, B, D, E, F, G, H, P, Q, R, S, Z, AB)
# are cyclically permuted:
#   row 16 <-> row 19            (2-cycle)
#   row 18 -> row 21 -> row 20 -> row 18   (3-cycle)
#
# Rather than physically moving cells (which could disturb styles/merges),
# we simply overwrite the per-row field values with their new contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) whose values move together with a record.
$cols = @("A","B","D","E","F","G","H","P","Q","R","S","Z","AB")

# Snapshot the current ("before") values for every row involved, keyed by
# row number and then column letter, so the writes below don't clobber
# data we still need to read.
$rowsInvolved = @(16, 18, 19, 20, 21)
$before = @{}
foreach ($r in $rowsInvolved) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value()
    }
    $before[$r] = $rowData
}

# Target row -> source row it should copy its field values from.
$mapping = @{ 16 = 19; 19 = 16; 18 = 20; 20 = 21; 21 = 18 }

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $srcData = $before[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $srcData[$c]
    }
}
